# Sync the "Fields" worksheet with the updated web-form config:
#  - remove the old "company" / "Company Title" field row
#  - give "company_type" a group_id
#  - add a new "new_field_2" field (select, with its own options) right after "company_type"
#  - replace the old "d-u-n-s_number" field with a new "new_field_3" text field
#  - refresh the selection / column width cosmetics left behind by the edit

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fields")

# 1. Drop the "company" / "Company Title" record - every row below shifts up by one.
$ws.Rows.Item(2).Delete()

# 2. "company_type" (now row 4) moves into group_id 1.
$ws.Cells.Item(4, 8).Value = 1

# 3. Insert the new "new_field_2" select field right after "company_type".
$ws.Rows.Item(5).Insert()
$ws.Cells.Item(5, 1).Value = "new_field_2"
$ws.Cells.Item(5, 2).Value = "This is a new field"
$ws.Cells.Item(5, 3).Value = "Yes"
$ws.Cells.Item(5, 4).Value = "select"
$ws.Cells.Item(5, 5).Value = "STRING"
$ws.Cells.Item(5, 6).Value = "Option 1, Option 2"
$ws.Cells.Item(5, 7).Value = 3
$ws.Cells.Item(5, 8).Value = 1

# 4. Replace "d-u-n-s_number" (now row 6) with the new "new_field_3" text field.
$ws.Cells.Item(6, 1).Value = "new_field_3"
$ws.Cells.Item(6, 2).Value = "Label 2"
$ws.Cells.Item(6, 3).Value = "No"
$ws.Cells.Item(6, 4).Value = "text"
$ws.Cells.Item(6, 5).Value = "STRING"
$ws.Cells.Item(6, 6).ClearContents()
$ws.Cells.Item(6, 7).Value = 4
$ws.Cells.Item(6, 8).ClearContents()

# 5. Column B needs to widen to fit the longer "This is a new field" label.
$ws.Columns.Item(2).ColumnWidth = 15.75

# 6. Leave the selection where the editor ended up when done.
$ws.Activate() | Out-Null
$ws.Range("I6").Select() | Out-Null
